$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Locations sheet: add new location rows (175-181) and extend the
# sheet with blank formula rows down to row 197.
# ---------------------------------------------------------------
$locations = $wb.Worksheets.Item("Locations")

$locRows = @(
    @{ Row=175; AreaId=10; WorldId=14; Name="HalloweenWrapping";   Display="Christmas Town: The Wrapping Room";  Hide=1 },
    @{ Row=176; AreaId=7;  WorldId=14; Name="HalloweenPlaza";      Display="Christmas Town: Christmas Tree Plaza"; Hide=1 },
    @{ Row=177; AreaId=15; WorldId=7;  Name="AgrabahShop2";        Display="The Peddler’s Shop";                 Hide=0 },
    @{ Row=178; AreaId=14; WorldId=7;  Name="AgrabahSand";         Display="Sandswept Ruins";                    Hide=0 },
    @{ Row=179; AreaId=11; WorldId=7;  Name="AgrabahRuin";         Display="Ruined Chamber";                     Hide=0 },
    @{ Row=180; AreaId=5;  WorldId=7;  Name="AgrabahAbove";        Display="Above the City";                     Hide=0 },
    @{ Row=181; AreaId=15; WorldId=10; Name="PrideSavannahBattle"; Display="The Savannah";                       Hide=0 }
)

foreach ($r in $locRows) {
    $row = $r.Row
    $locations.Cells.Item($row, 2).Value = $r.AreaId
    $locations.Cells.Item($row, 3).Value = $r.WorldId
    $locations.Cells.Item($row, 4).Value = $r.Name
    $locations.Cells.Item($row, 5).Value = $r.Display
    $locations.Cells.Item($row, 6).Value = $r.Hide
    $locations.Cells.Item($row, 1).Formula = "=B$row+C$row*256"

    $gFormula = '=_xlfn.CONCAT( ,A' + $row + ',": { ""worldId"": ",C' + $row + ',", ""name"": """,D' + $row + ',""", ""display"": """,E' + $row + ',""", ""areaId"": ",B' + $row + ',", ""hideWorld"": ",F' + $row + ',", },")'
    $locations.Cells.Item($row, 7).Formula = $gFormula
}

for ($row = 182; $row -le 197; $row++) {
    $locations.Cells.Item($row, 1).Formula = "=B$row+C$row*256"
}

# ---------------------------------------------------------------
# Events sheet: add new achievement/event rows (25-27).
# ---------------------------------------------------------------
$events = $wb.Worksheets.Item("Events")

$events.Cells.Item(25, 1).Value = 40          # text "40" -> "0x40"
$events.Cells.Item(25, 2).Value = "HalloweenPlaza"
$events.Cells.Item(25, 3).Value = "Experiment"
$events.Cells.Item(25, 4).Value = 1
$events.Cells.Item(25, 5).Value = 25
$events.Cells.Item(25, 6).Value = 1

$events.Cells.Item(26, 1).Value = "3d"
$events.Cells.Item(26, 2).Value = "AgrabahSand"
$events.Cells.Item(26, 3).Value = "Sandswept Ruins Escape"
$events.Cells.Item(26, 4).Value = 0
$events.Cells.Item(26, 5).Value = 0
$events.Cells.Item(26, 6).Value = 0

$events.Cells.Item(27, 1).Value = "3e"
$events.Cells.Item(27, 2).Value = "AgrabahAbove"
$events.Cells.Item(27, 3).Value = "Jafar"
$events.Cells.Item(27, 4).Value = 1
$events.Cells.Item(27, 5).Value = 10
$events.Cells.Item(27, 6).Value = 1

foreach ($row in 25..27) {
    $gFormula = '=_xlfn.CONCAT("""",C' + $row + ',""": { ""locationCode"": """,B' + $row + ',""", ""eventId"": 0x",A' + $row + ',", ""score"": ",E' + $row + ',", ""isBoss"": ",D' + $row + ',", ""lv1"": ",F' + $row + ',' + '" }, ")'
    $events.Cells.Item($row, 7).Formula = $gFormula
}

# ---------------------------------------------------------------
# Restore view state (selection) for the Locations sheet.
# ---------------------------------------------------------------
$locations.Range("F184").Select()

Write-Host "Done"
